# Apply updated crypto price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.002.86"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "2.498.22"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'321.03"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "'107.84"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").Value = "'0.525"
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.539"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "'39.79"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("E11").Value = "  +8.13%  "
$ws.Range("D12").Value = "'0.0815"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "'7.15"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "2.890.09"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "2.504.65"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "'0.842"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "47.937.00"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").Value = "'13.04"
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("D20").Value = "'6.75"
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").Value = "0.0₃0943"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'71.89"
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'276.51"
$ws.Range("E24").Value = "  +11.89%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'25.70"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.28"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'9.79"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.139"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'35.23"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("D32").Value = "'49.15"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").Value = "'19.49"
$ws.Range("E33").Value = "  -4.85%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'1.01"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'5.33"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").Value = "'0.0783"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").Value = "'4.63"
$ws.Range("E38").Value = "  -3.50%  "
$ws.Range("D39").Value = "'2.92"
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").Value = "'121.16"
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").Value = "'21.49"
$ws.Range("E43").Value = "  -6.61%  "
$ws.Range("D44").Value = "'0.0303"
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").Value = "1.998.85"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'3.15"
$ws.Range("E46").Value = "  +4.23%  "
$ws.Range("D47").Value = "'1.85"
$ws.Range("E47").Value = "  +3.24%  "
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("D49").Value = "'8.94"
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").Value = "'5.16"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").Value = "'79.77"
$ws.Range("E51").Value = "  +2.83%  "
